$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 88.94787128385633
$ws.Range("B2").Value = 3.2672515866804757
$ws.Range("C2").Value = -19.80112294042058
$ws.Range("D2").Value = 6.1246374422336958
$ws.Range("E2").Value = 331.1098736231051
$ws.Range("F2").Value = 7.785145470917735
$ws.Range("G2").Value = 540.67252418319038
$ws.Range("H2").Value = 8.6483430239752437
$ws.Range("I2").Value = 0.22287795462860235
$ws.Range("J2").Value = 0.22287795462860235
$ws.Range("K2").Value = 7.8821297219057387
$ws.Range("L2").Value = 1.1260185317008198
$ws.Range("M2").Value = 0.98780601013284786
$ws.Range("N2").Value = -0.63612651092666672
